# aggiornamento fino a 20/09/2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows to append (row index, date serial, B, C, D)
$rows = @(
    @(375, 44449, 0, 1, 40.79967360261118),
    @(376, 44450, 0, 1, 40.79967360261118),
    @(377, 44451, 0, 1, 40.79967360261118),
    @(378, 44452, 0, 0, 0),
    @(379, 44453, 0, 0, 0),
    @(380, 44454, 0, 0, 0),
    @(381, 44455, 0, 0, 0),
    @(382, 44456, 0, 0, 0),
    @(383, 44457, 0, 0, 0),
    @(384, 44458, 2, 2, 81.59934720522236),
    @(385, 44459, 0, 2, 81.59934720522236)
)

$styleSource = $ws.Cells.Item(374, 1)

foreach ($r in $rows) {
    $rowIndex = $r[0]
    $dateVal = $r[1]
    $bVal = $r[2]
    $cVal = $r[3]
    $dVal = $r[4]

    $aCell = $ws.Cells.Item($rowIndex, 1)
    $styleSource.Copy($aCell)
    $aCell.Value = $dateVal

    $ws.Cells.Item($rowIndex, 2).Value = $bVal
    $ws.Cells.Item($rowIndex, 3).Value = $cVal
    $ws.Cells.Item($rowIndex, 4).Value = $dVal
}
